$wb = $excel.ActiveWorkbook

# --- RoundTrip sheet (new flight search data) ---
$ws1 = $wb.Worksheets.Item("RoundTrip")

# Row 2: New Delhi -> Patna
$ws1.Cells.Item(2, 2).Value = "Patna"
$ws1.Cells.Item(2, 3).Value = "01/09/2020"
$ws1.Cells.Item(2, 4).Value = "01/10/2020"
$ws1.Cells.Item(2, 8).Value = "No"

# Row 3: Mumbai -> Bhopal
$ws1.Cells.Item(3, 3).Value = "18/09/2020"
$ws1.Cells.Item(3, 4).Value = "17/11/2020"
$ws1.Cells.Item(3, 8).Value = "No"

# Row 4: Hyderabad -> Visakhapatnam
$ws1.Cells.Item(4, 3).Value = "09/11/2020"
$ws1.Cells.Item(4, 4).Value = "04/12/2020"

# Row 5: Bangalore -> Kolkata
$ws1.Cells.Item(5, 3).Value = "04/10/2020"
$ws1.Cells.Item(5, 4).Value = "01/01/2021"
$ws1.Cells.Item(5, 8).NumberFormat = "@"
$ws1.Cells.Item(5, 8).Value = "Yes"

# --- HotelSearch sheet (logout -> executed as Yes now) ---
$ws2 = $wb.Worksheets.Item("HotelSearch")

$ws2.Cells.Item(2, 2).Value = "15/08/2020"
$ws2.Cells.Item(2, 3).Value = "02/09/2020"

$ws2.Cells.Item(3, 2).Value = "25/10/2020"
$ws2.Cells.Item(3, 3).Value = "05/11/2020"
$ws2.Cells.Item(3, 5).Value = "Yes"

$ws2.Cells.Item(4, 5).Value = "Yes"

# --- Update selections to match the edited cells ---
$ws2.Range("C9").Select()

# RoundTrip becomes the active sheet/tab, with F9 selected
$ws1.Activate()
$ws1.Range("F9").Select()
